$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New total extent is A1:D115 (2 new days of data appended at the front of
# this "missing Jan 5th 2021" stretch, and the whole 7-day rolling window
# recomputed / shifted by one day, plus a brand new last row).

# Make sure the two brand-new rows (114 and 115) inherit the same cell
# style (date format, border, bold-less body font) as the rest of column A
# by copying an existing dated cell's formatting down first.
$ws.Range("A113").Copy($ws.Range("A114"))
$ws.Range("A113").Copy($ws.Range("A115"))

# r, date(serial), nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila ab.
$rows = @(
  @(93, 44235, 0, 1, 28.87669650591972),
  @(94, 44236, 0, 2, 57.75339301183945),
  @(95, 44237, 0, 2, 57.75339301183945),
  @(96, 44238, 0, 1, 28.87669650591972),
  @(97, 44239, 1, 2, 57.75339301183945),
  @(98, 44240, 0, 2, 57.75339301183945),
  @(99, 44241, 0, 2, 57.75339301183945),
  @(100, 44242, 1, 3, 86.63008951775916),
  @(101, 44243, 0, 3, 86.63008951775916),
  @(102, 44244, 0, 3, 86.63008951775916),
  @(103, 44245, 1, 3, 86.63008951775916),
  @(104, 44246, 1, 3, 86.63008951775916),
  @(105, 44247, 0, 3, 86.63008951775916),
  @(106, 44248, 0, 3, 86.63008951775916),
  @(107, 44249, 1, 2, 57.75339301183945),
  @(108, 44250, 0, 1, 28.87669650591972),
  @(109, 44251, 0, 1, 28.87669650591972),
  @(110, 44252, 0, 2, 57.75339301183945),
  @(111, 44253, 0, 2, 57.75339301183945),
  @(112, 44254, 0, 2, 57.75339301183945)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
}

# Rows 113-115 have no 7-day rolling sum yet (trailing edge of the series),
# so only the date and the "nuovi pos." counter move.
$tailRows = @(
  @(113, 44255, 1),
  @(114, 44256, 1),
  @(115, 44257, 0)
)

foreach ($r in $tailRows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
}
